$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Row 2 updates
# ---------------------------------------------------------------------------
$ws.Cells.Item(2, 2).Value = 8                      # B2
$ws.Cells.Item(2, 8).Value = ""                     # H2 (was 5, now blank)
$ws.Cells.Item(2, 9).Value = "2022-08-27T17:46:58.290699Z"   # I2
$ws.Cells.Item(2, 10).Value = "2022-08-27T17:46:58.290744Z" # J2

# ---------------------------------------------------------------------------
# Row 3 updates
# ---------------------------------------------------------------------------
$ws.Cells.Item(3, 2).Value = 9                      # B3
$ws.Cells.Item(3, 3).Value = "test1"                # C3
$ws.Cells.Item(3, 4).Value = ""                     # D3 (was abc0100002)
$ws.Cells.Item(3, 8).Value = 7                      # H3
$ws.Cells.Item(3, 9).Value = "2022-08-27T17:47:47.062414Z"   # I3
$ws.Cells.Item(3, 10).Value = "2022-08-29T22:25:59.011370Z"  # J3
$ws.Cells.Item(3, 71).Value = "0"                   # BS3
$ws.Cells.Item(3, 72).Value = "test3"               # BT3

# ---------------------------------------------------------------------------
# Row 4 updates
# ---------------------------------------------------------------------------
$ws.Cells.Item(4, 2).Value = 10                     # B4
$ws.Cells.Item(4, 3).Value = "test2"                # C4
$ws.Cells.Item(4, 8).Value = 2                      # H4
$ws.Cells.Item(4, 9).Value = "2022-08-27T17:48:10.571237Z"   # I4
$ws.Cells.Item(4, 10).Value = "2022-08-29T22:25:52.385998Z"  # J4
$ws.Cells.Item(4, 70).Value = "test1"                # BR4
$ws.Cells.Item(4, 71).Value = "sadda"                # BS4
$ws.Cells.Item(4, 72).Value = "test2"                # BT4

# ---------------------------------------------------------------------------
# New row 5
# ---------------------------------------------------------------------------
$ws.Cells.Item(5, 1).Value = 3
$ws.Cells.Item(5, 1).Style = $ws.Cells.Item(4, 1).Style
$ws.Cells.Item(5, 2).Value = 11
$ws.Cells.Item(5, 3).Value = "test3"
$ws.Cells.Item(5, 8).Value = 2
$ws.Cells.Item(5, 9).Value = "2022-08-27T17:51:40.745322Z"
$ws.Cells.Item(5, 10).Value = "2022-08-29T22:25:45.884907Z"
$ws.Cells.Item(5, 21).Value = $false               # U5
$ws.Cells.Item(5, 24).Value = $true                 # X5
$ws.Cells.Item(5, 69).Value = 1                     # BQ5
$ws.Cells.Item(5, 70).Value = "test2"               # BR5
$ws.Cells.Item(5, 71).Value = "asdad"                # BS5
$ws.Cells.Item(5, 72).Value = "test"                # BT5

# ---------------------------------------------------------------------------
# New row 6
# ---------------------------------------------------------------------------
$ws.Cells.Item(6, 1).Value = 4
$ws.Cells.Item(6, 1).Style = $ws.Cells.Item(4, 1).Style
$ws.Cells.Item(6, 2).Value = 12
$ws.Cells.Item(6, 3).Value = "mamad"
$ws.Cells.Item(6, 21).Value = $false               # U6
$ws.Cells.Item(6, 24).Value = $true                 # X6
$ws.Cells.Item(6, 9).Value = "2022-08-27T17:58:21.421881Z"
$ws.Cells.Item(6, 10).Value = "2022-08-29T22:25:37.756867Z"
$ws.Cells.Item(6, 69).Value = 1                     # BQ6
$ws.Cells.Item(6, 70).Value = "test3"               # BR6
$ws.Cells.Item(6, 71).Value = "0"                   # BS6
$ws.Cells.Item(6, 72).Value = "admin"               # BT6

Write-Output "facility seg view rows updated"
